# Slide 9, shape 6 is the "Picture 8" image (p:pic) that gets resized/
# repositioned and given a red outline in the target revision.
$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(9)
$sh = $s.Shapes.Item(6)

# Reposition / resize the picture (point values chosen so that, after the
# host's internal EMU conversion, they land on the exact target EMU values:
# off x=1115616 y=2083666, ext cx=6912768 cy=3626241).
$sh.Left   = 87.84381866455078
$sh.Top    = 164.06822204589844
$sh.Width  = 544.3125
$sh.Height = 285.53082275390625

# Give the picture a solid red (C00000) outline.
$sh.Line.Visible = $true
$sh.Line.ForeColor.RGB = 192
